# Fruta / hortaliza, semanal
# A new daily price record was added for "Vega Modelo de Temuco" / Plátano.
# It is inserted as row 247 (pushing the former rows 247-337 down to 248-338),
# which is why the sheet dimension grows from A1:T337 to A1:T338.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 247, shifting existing rows 247-337 down to 248-338.
$ws.Rows.Item(247).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(247, 1).Value = 10
$ws.Cells.Item(247, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(247, 3).Value = 'La Araucanía'
$ws.Cells.Item(247, 4).Value = 44468
$ws.Cells.Item(247, 5).Value = 9
$ws.Cells.Item(247, 6).Value = 'Fruta'
$ws.Cells.Item(247, 7).Value = 100108
$ws.Cells.Item(247, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(247, 9).Value = 100108006
$ws.Cells.Item(247, 10).Value = 'Plátano'
$ws.Cells.Item(247, 11).Value = 'Sin especificar'
$ws.Cells.Item(247, 12).Value = 'Pintón'
$ws.Cells.Item(247, 13).Value = 350
$ws.Cells.Item(247, 14).Value = 15000
$ws.Cells.Item(247, 15).Value = 16000
$ws.Cells.Item(247, 16).Value = 15429
$ws.Cells.Item(247, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(247, 18).Value = 'Ecuador'
$ws.Cells.Item(247, 19).Value = 771
$ws.Cells.Item(247, 20).Value = 20
